{"js": "// Apply the Mastercard Orange (#FF5F00) heading template update to the\n// Heading 1 - Heading 4 paragraph styles:\n//   - font becomes \"Calibri Light\" (ascii + hAnsi; theme fallbacks kept)\n//   - color becomes a hard-coded RGB #FF5F00 (theme color reference dropped)\n//   - font size bumps up (Heading1: 18pt, Heading2: 15pt, Heading3: 13pt, Heading4: 12pt)\nconst styles = context.document.getStyles();\n\nconst headingSizes = {\n  \"Heading 1\": 18,\n  \"Heading 2\": 15,\n  \"Heading 3\": 13,\n  \"Heading 4\": 12,\n};\n\nfor (const [name, size] of Object.entries(headingSizes)) {\n  const style = styles.getByNameOrNullObject(name);\n  style.font.name = \"Calibri Light\";\n  style.font.color = \"#FF5F00\";\n  style.font.size = size;\n}\n\nawait context.sync();\n", "ps1": "# Apply the Mastercard Orange (#FF5F00) heading template update to\n# Heading 1 - Heading 4 paragraph styles:\n#   - font becomes \"Calibri Light\" (ascii + hAnsi; theme fallbacks kept)\n#   - color becomes a hard-coded RGB FF5F00 (theme color reference dropped)\n#   - font size bumps up (Heading1: 18pt, Heading2: 15pt, Heading3: 13pt, Heading4: 12pt)\n$d = $word.ActiveDocument\n\n# Word's Font.Color is stored BGR (0x00BBGGRR) -- FF5F00 (R=FF,G=5F,B=00) -> 0x00005FFF\n$mastercardOrange = 0x00005FFF\n\n$headingSizes = @{\n    \"Heading 1\" = 18\n    \"Heading 2\" = 15\n    \"Heading 3\" = 13\n    \"Heading 4\" = 12\n}\n\nforeach ($name in $headingSizes.Keys) {\n    $style = $d.Styles($name)\n    $style.Font.Name = \"Calibri Light\"\n    $style.Font.Color = $mastercardOrange\n    $style.Font.Size = $headingSizes[$name]\n}\n"}
